$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 13 new rows right above the current row 2 (pushes old data rows down).
$ws.Rows("2:14").Insert()

# The insert picks up the header row's bold formatting by default; strip it
# back to the unstyled look the rest of the numeric data uses.
$ws.Range("A2:C14").ClearFormats()

# New data for the freshly inserted rows 2-14.
$newData = @(
    @(-0.1624901592731475, 0.1201877370476722, 0.1788308024406433),
    @(0.06856962293386459, -0.06704246252775189, 0.1534798890352249),
    @(0.0235183127224445, -0.0152716310694813, 0.030695978552103),
    @(0.078801617026329, 0.009010262787342, -0.0384845100343227),
    @(0.0911716371774673, 0.0236710291355848, -0.0128281703218817),
    @(0.0241291765123605, -0.0087048299610614, -0.0103847095742821),
    @(0.0088575463742017, 0.0404698215425014, 0.0618501044809818),
    @(-0.0404698215425014, -0.0221438650041818, 0.0143553335219621),
    @(0.0478002056479454, -0.022754730656743, 0.08124507963657369),
    @(0.0325285755097866, 0.0065668015740811, 0.0574213340878486),
    @(-0.0219911485910415, 0.0610865242779254, 0.0355829000473022),
    @(-0.0038179077673703, 0.0641408488154411, -0.0708603709936142),
    @(0.0064140851609408, -0.0274889357388019, -0.0403171069920063)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = 2 + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# The old trailing 3 rows (previously rows 19-21) are no longer part of the
# dataset after the shift - they now sit at rows 32-34, so remove them.
$ws.Rows("32:34").Delete()
